$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.789.84"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "2.311.41"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.43"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.505"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.98"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0783"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("E13").Value = "  -0.01%  "
$ws.Range("E14").Value = "  -1.70%  "
$ws.Range("D15").Value = "2.671.54"
$ws.Range("D16").Value = "2.304.08"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("E17").Value = "  +1.99%  "
$ws.Range("D18").Value = "42.730.44"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("E23").Value = "  +6.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("E26").Value = "  +1.63%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.19%  "
$ws.Range("E28").Value = "  +15.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "166.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.13"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.11"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.97%  "
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("E35").Value = "  +0.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0698"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.01%  "
$ws.Range("E37").Value = "  -0.87%  "
$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.04%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.100"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +1.01%  "
$ws.Range("E41").Value = "  -0.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.91"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +15.02%  "
$ws.Range("D43").Value = "1.926.16"
$ws.Range("E43").Value = "  -3.33%  "
$ws.Range("E44").Value = "  -0.23%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("E46").Value = "  -2.75%  "
$ws.Range("E47").Value = "  -0.86%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.89"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.48%  "
$ws.Range("D49").Value = "2.540.51"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.34"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.02%  "
